$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 149339.7
$ws.Range("I6").Value = 211672.14
$ws.Range("J6").Value = 3897.3333
$ws.Range("K6").Value = 635016.42
$ws.Range("L6").Value = 11691.9999
$ws.Range("M6").Value = -634904.42
$ws.Range("N6").Value = -11915.9999

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 293.75
$ws.Range("I9").Value = 321.42856
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 321.42856
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -152.42856
$ws.Range("N9").Value = -438

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 16318.903
$ws.Range("I12").Value = 177.36
$ws.Range("J12").Value = 83575.336
$ws.Range("K12").Value = 177.36
$ws.Range("L12").Value = 83575.336
$ws.Range("M12").Value = -7.360000000000014

# ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 8222.182000000001
$ws.Range("I21").Value = 8472
$ws.Range("J21").Value = 8166.6665
$ws.Range("K21").Value = 8472
$ws.Range("L21").Value = 8166.6665
$ws.Range("M21").Value = -8004
$ws.Range("N21").Value = -9102.666499999999

# ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 8222.182000000001
$ws.Range("I23").Value = 8472
$ws.Range("J23").Value = 8166.6665
$ws.Range("K23").Value = 8472
$ws.Range("L23").Value = 8166.6665
$ws.Range("M23").Value = -8238
$ws.Range("N23").Value = -8634.666499999999

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1792461.1
$ws.Range("I38").Value = 2932685.8
$ws.Range("J38").Value = 679.4286
$ws.Range("K38").Value = 8798057.399999999
$ws.Range("L38").Value = 2038.2858
$ws.Range("M38").Value = -8797685.399999999
$ws.Range("N38").Value = -2782.2858

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 722732.1
$ws.Range("I58").Value = 1262906.2
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 3788718.6
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -3788568.6
$ws.Range("N58").Value = -7800

# ALC row 81
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 36500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 36500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 36500
$ws.Range("N81").Value = -38496

# ALC row 84
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H84").Value = 36500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 36500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 109500
$ws.Range("N84").Value = -119484

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 31034.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 31034.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 31034.5
$ws.Range("N87").Value = -33530.5

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 31034.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 31034.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 93103.5
$ws.Range("N90").Value = -105583.5

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2341.8298
$ws.Range("I45").Value = 1974.8064
$ws.Range("J45").Value = 3052.9375
$ws.Range("K45").Value = 1974.8064
$ws.Range("L45").Value = 3052.9375
$ws.Range("M45").Value = -1597.8064
$ws.Range("N45").Value = -3806.9375

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1238.7142
$ws.Range("I99").Value = 823.3333
$ws.Range("J99").Value = 1986.4
$ws.Range("K99").Value = 823.3333
$ws.Range("L99").Value = 1986.4
$ws.Range("M99").Value = 674.6667
$ws.Range("N99").Value = -4982.4

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 41670304
$ws.Range("I107").Value = 55557624
$ws.Range("J107").Value = 8332.5
$ws.Range("K107").Value = 55557624
$ws.Range("L107").Value = 8332.5
$ws.Range("M107").Value = -55555704
$ws.Range("N107").Value = -12172.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17401.484
$ws.Range("I31").Value = 34550.7
$ws.Range("J31").Value = 3110.4722
$ws.Range("K31").Value = 34550.7
$ws.Range("L31").Value = 3110.4722
$ws.Range("M31").Value = -34255.7
$ws.Range("N31").Value = -3700.4722

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 17401.484
$ws.Range("I34").Value = 34550.7
$ws.Range("J34").Value = 3110.4722
$ws.Range("K34").Value = 34550.7
$ws.Range("L34").Value = 3110.4722
$ws.Range("M34").Value = -34348.7
$ws.Range("N34").Value = -3514.4722

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3448.1943
$ws.Range("I107").Value = 6316.5557
$ws.Range("J107").Value = 579.8333
$ws.Range("K107").Value = 6316.5557
$ws.Range("L107").Value = 579.8333
$ws.Range("M107").Value = -4396.5557
$ws.Range("N107").Value = -4419.8333

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2625.2778
$ws.Range("I132").Value = 2585.6428
$ws.Range("J132").Value = 2764
$ws.Range("K132").Value = 7756.928400000001
$ws.Range("L132").Value = 8292
$ws.Range("M132").Value = -5226.928400000001
$ws.Range("N132").Value = -13352

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1719.9333
$ws.Range("I134").Value = 778.6
$ws.Range("J134").Value = 3602.6
$ws.Range("K134").Value = 2335.8
$ws.Range("L134").Value = 10807.8
$ws.Range("M134").Value = 199.1999999999998
$ws.Range("N134").Value = -15877.8

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 962
$ws.Range("I63").Value = 962
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2886
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2137

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 962
$ws.Range("I66").Value = 962
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8658
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4914

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2327.3147
$ws.Range("I68").Value = 1761.762
$ws.Range("J68").Value = 2501.9707
$ws.Range("K68").Value = 5285.286
$ws.Range("L68").Value = 7505.9121
$ws.Range("M68").Value = -4474.286
$ws.Range("N68").Value = -9127.9121

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2327.3147
$ws.Range("I71").Value = 1761.762
$ws.Range("J71").Value = 2501.9707
$ws.Range("K71").Value = 15855.858
$ws.Range("L71").Value = 22517.7363
$ws.Range("M71").Value = -11799.858
$ws.Range("N71").Value = -30629.7363

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 321.15384
$ws.Range("I2").Value = 366
$ws.Range("J2").Value = 74.5
$ws.Range("K2").Value = 366
$ws.Range("L2").Value = 74.5
$ws.Range("M2").Value = -253
$ws.Range("N2").Value = -300.5

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3306.3333
$ws.Range("I43").Value = 1933.3334
$ws.Range("J43").Value = 4679.3335
$ws.Range("K43").Value = 1933.3334
$ws.Range("L43").Value = 4679.3335
$ws.Range("M43").Value = -1782.3334
$ws.Range("N43").Value = -4981.3335

# GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14839.8
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 14839.8
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 14839.8
$ws.Range("N46").Value = -15151.8

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18899.75
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 18899.75
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 18899.75
$ws.Range("N57").Value = -20539.75

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2507.7273
$ws.Range("I80").Value = 2648.125
$ws.Range("J80").Value = 2133.3333
$ws.Range("K80").Value = 2648.125
$ws.Range("L80").Value = 2133.3333
$ws.Range("M80").Value = -1650.125
$ws.Range("N80").Value = -4129.3333

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2507.7273
$ws.Range("I83").Value = 2648.125
$ws.Range("J83").Value = 2133.3333
$ws.Range("K83").Value = 13240.625
$ws.Range("L83").Value = 10666.6665
$ws.Range("M83").Value = -8248.625
$ws.Range("N83").Value = -20650.6665

# LTW row 9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5915.7144
$ws.Range("I9").Value = 177.5
$ws.Range("J9").Value = 13566.667
$ws.Range("K9").Value = 177.5
$ws.Range("L9").Value = 13566.667
$ws.Range("M9").Value = 46.5
$ws.Range("N9").Value = -14014.667

Write-Output "edits applied"
